$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from column F (the old column D, now shifted) into the new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new columns D (latest quarter) and E (prior quarter) with their reported values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 859500
$ws.Range("E8").Value = 990300
$ws.Range("D9").Value = 630300
$ws.Range("E9").Value = 749000
$ws.Range("D10").Value = 229200
$ws.Range("E10").Value = 241300
$ws.Range("D12").Value = 'NA'
$ws.Range("E12").Value = 'NA'
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 400
$ws.Range("E14").Value = 1500
$ws.Range("D15").Value = 14100
$ws.Range("E15").Value = 13800
$ws.Range("D17").Value = 818800
$ws.Range("E17").Value = 940500
$ws.Range("D18").Value = 40700
$ws.Range("E18").Value = 49800
$ws.Range("D20").Value = 2800
$ws.Range("E20").Value = 2900
$ws.Range("D21").Value = 60400
$ws.Range("E21").Value = 69400
$ws.Range("D22").Value = 6100
$ws.Range("E22").Value = 5900
$ws.Range("D23").Value = 37500
$ws.Range("E23").Value = 46800
$ws.Range("D24").Value = 9800
$ws.Range("E24").Value = 11000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 27600
$ws.Range("E26").Value = 35900
$ws.Range("D27").Value = 27600
$ws.Range("E27").Value = 35900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 500
$ws.Range("E29").Value = 'NA'
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2800
$ws.Range("E32").Value = -2900
$ws.Range("D33").Value = 28100
$ws.Range("E33").Value = 35900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 28100
$ws.Range("E35").Value = 35900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 150700
$ws.Range("E41").Value = 57700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 330800
$ws.Range("E43").Value = 413600
$ws.Range("D44").Value = 309300
$ws.Range("E44").Value = 356600
$ws.Range("D45").Value = 56200
$ws.Range("E45").Value = 62300
$ws.Range("D46").Value = 847000
$ws.Range("E46").Value = 890200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 294300
$ws.Range("E48").Value = 297400
$ws.Range("D49").Value = 421900
$ws.Range("E49").Value = 427700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 12900
$ws.Range("E52").Value = 11500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1576100
$ws.Range("E54").Value = 1626800
$ws.Range("D57").Value = 123500
$ws.Range("E57").Value = 205300
$ws.Range("D58").Value = 6700
$ws.Range("E58").Value = 7000
$ws.Range("D59").Value = 166000
$ws.Range("E59").Value = 161300
$ws.Range("D60").Value = 296200
$ws.Range("E60").Value = 373600
$ws.Range("D61").Value = 354000
$ws.Range("E61").Value = 355600
$ws.Range("D62").Value = 51200
$ws.Range("E62").Value = 50900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 701500
$ws.Range("E66").Value = 780200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 210300
$ws.Range("E72").Value = 182200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 874700
$ws.Range("E76").Value = 846700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 28100
$ws.Range("E81").Value = 35900
$ws.Range("D83").Value = 16800
$ws.Range("E83").Value = 16600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 99400
$ws.Range("E89").Value = 59700
$ws.Range("D91").Value = -12500
$ws.Range("E91").Value = -16400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -4200
$ws.Range("E94").Value = -12200
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2100
$ws.Range("E100").Value = -4200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 93000
$ws.Range("E102").Value = 43300
$ws.Range("F43").Value = 423100
$ws.Range("G43").Value = 392700
$ws.Range("F45").Value = 72200
$ws.Range("G45").Value = 55100
